# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates the G column ("K") values for rows 2-37 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 5
    3  = 8
    4  = 4
    5  = 6
    6  = 13
    7  = 8
    8  = 5
    9  = 5
    10 = 7
    11 = 4
    12 = 14
    13 = 3
    14 = 8
    15 = 7
    16 = 6
    17 = 4
    18 = 12
    19 = 4
    20 = 7
    21 = 7
    22 = 5
    23 = 6
    24 = 6
    25 = 4
    26 = 7
    27 = 3
    28 = 4
    29 = 5
    30 = 7
    31 = 6
    32 = 4
    33 = 6
    34 = 1
    35 = 6
    36 = 1
    37 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
